$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-09-03 Wednesday"; new="2025-09-04 Thursday"},
    @{old="817×8="; new="692×9="},
    @{old="779×9="; new="630×5="},
    @{old="936×8="; new="553×8="},
    @{old="220×6="; new="980×3="},
    @{old="984×6="; new="871×3="},
    @{old="564×2="; new="948×2="},
    @{old="944×2="; new="168×3="},
    @{old="555×5="; new="843×2="},
    @{old="357×6="; new="675×5="},
    @{old="982×4="; new="962×6="},
    @{old="595×6="; new="984×2="},
    @{old="931×5="; new="682×9="},
    @{old="285×7="; new="347×7="},
    @{old="816×3="; new="234×9="},
    @{old="755×5="; new="455×9="},
    @{old="825×4="; new="933×8="},
    @{old="494×5="; new="565×5="},
    @{old="633×4="; new="853×5="},
    @{old="553×3="; new="611×5="},
    @{old="261×3="; new="782×7="},
    @{old="850×2="; new="920×6="},
    @{old="951×8="; new="601×8="},
    @{old="476×9="; new="791×2="},
    @{old="492×9="; new="844×4="},
    @{old="738×9="; new="699×2="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
